# Update "想去人数" (want-to-go count) values in column F across the
# 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# ---- 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 16
$ws1.Range("F5").Value = 68
$ws1.Range("F6").Value = 4
$ws1.Range("F7").Value = 1303
$ws1.Range("F8").Value = 510
$ws1.Range("F10").Value = 1284
$ws1.Range("F13").Value = 1089
$ws1.Range("F14").Value = 14
$ws1.Range("F17").Value = 93
$ws1.Range("F18").Value = 243
$ws1.Range("F19").Value = 1655
$ws1.Range("F21").Value = 266
$ws1.Range("F22").Value = 204
$ws1.Range("F23").Value = 1997
$ws1.Range("F26").Value = 920
$ws1.Range("F27").Value = 1206
$ws1.Range("F30").Value = 2809
$ws1.Range("F31").Value = 1603
$ws1.Range("F32").Value = 81
$ws1.Range("F33").Value = 113
$ws1.Range("F34").Value = 647
$ws1.Range("F36").Value = 1789
$ws1.Range("F37").Value = 887
$ws1.Range("F38").Value = 1802
$ws1.Range("F39").Value = 200
$ws1.Range("F42").Value = 38
$ws1.Range("F43").Value = 846
$ws1.Range("F44").Value = 788
$ws1.Range("F45").Value = 1000
$ws1.Range("F46").Value = 56
$ws1.Range("F47").Value = 434
$ws1.Range("F48").Value = 3325

# ---- 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 187
$ws2.Range("F12").Value = 793
$ws2.Range("F20").Value = 29
$ws2.Range("F21").Value = 18

# ---- 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 16
$ws4.Range("F5").Value = 68
$ws4.Range("F7").Value = 187
$ws4.Range("F8").Value = 1303
$ws4.Range("F9").Value = 510
$ws4.Range("F11").Value = 1284
$ws4.Range("F14").Value = 1089
$ws4.Range("F15").Value = 14
$ws4.Range("F18").Value = 93
$ws4.Range("F20").Value = 243
$ws4.Range("F21").Value = 1655
$ws4.Range("F23").Value = 266
$ws4.Range("F24").Value = 204
$ws4.Range("F25").Value = 1997
$ws4.Range("F28").Value = 1206
$ws4.Range("F29").Value = 2809
$ws4.Range("F30").Value = 1603
$ws4.Range("F31").Value = 81
$ws4.Range("F32").Value = 113
$ws4.Range("F33").Value = 793
$ws4.Range("F35").Value = 647
$ws4.Range("F37").Value = 1789
$ws4.Range("F39").Value = 887
$ws4.Range("F40").Value = 1802
$ws4.Range("F42").Value = 846
$ws4.Range("F43").Value = 788
$ws4.Range("F44").Value = 1000
$ws4.Range("F45").Value = 434
$ws4.Range("F46").Value = 29
$ws4.Range("F47").Value = 18
$ws4.Range("F48").Value = 3325
